$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "'32297400"
$ws.Range("B3").Value = "'32297401"
$ws.Range("B4").Value = "'32297402"
